$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing "fully connected" rows
$ws.Range("C10").Value = 1000
$ws.Range("C11").Value = 300

# Add a new "fully connected" layer row (row 12)
$ws.Range("B12").Value = "fully connected"
$ws.Range("C12").Value = 200
$ws.Range("D12").Value = 1
$ws.Range("E12").Value = 1
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 1
$ws.Range("H12").Value = 1

# Add a new blank row (row 17) below the table, mirroring the existing blank rows
# (apply the same center-alignment formatting used by the other blank rows so the
# cells materialize without introducing stray values)
$ws.Range("C17:K17").HorizontalAlignment = -4108

# Move the active selection to the newly edited cell
$ws.Range("D12").Select()
